$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": porcelanato sales for CHONTASI SIMBAÑA SILVIA JANETH (row 7)
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M7").Value = 109.62
$wsVentasGrupo.Range("M23").Value = "6 de 21"

# Sheet "VENTA MENSUAL": septiembre sales for the same client (row 7) and column total (row 23)
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F7").Value = 109.62
$wsVentaMensual.Range("F23").Value = 23999.03

# Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO group row (12) and TOTAL row (15)
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D12").Value = 21905.78
$wsCumplimiento.Range("E12").Value = 14917.8630921171
$wsCumplimiento.Range("F12").Value = 0.5948835628566421

$wsCumplimiento.Range("D15").Value = 23999.03
$wsCumplimiento.Range("E15").Value = 31425.71316613378
$wsCumplimiento.Range("F15").Value = 0.4330020967000916
